$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style from an existing header cell (H1) to the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# New column data values
$I = @(7, 7, 6, 4, 7, 4, 7, 4, 8, 9, 5, 8, 7, 7, 7, 4)
$J = @(8, 7, 6, 5, 7, 5, 7, 5, 8, 9, 5, 8, 7, 7, 7, 4)

for ($i = 0; $i -lt $I.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $I[$i]
    $ws.Cells.Item($row, 10).Value = $J[$i]
}
